$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (append "1" to each existing value)
$ws.Range("B2").Value = "егорdsa2d1"
$ws.Range("C2").Value = "губин2d1"
$ws.Range("D2").Value = "выфывфы2d1"
$ws.Range("E2").Value = "выфв2d1"
$ws.Range("F2").Value = "выфвыфв1"
$ws.Range("G2").Value = "вфывф1"

# Cells that get a purely numeric-looking text value need to be forced to
# text format first so Excel doesn't silently convert them to numbers.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B5:G5").NumberFormat = "@"
$ws.Range("B6:G6").NumberFormat = "@"

# Update row 3: only B3 changes
$ws.Range("B3").Value = "3222"

# Update row 5 values: old "1" -> "36" (these are the values that used to live in row 6)
$ws.Range("B5").Value = "36"
$ws.Range("C5").Value = "36"
$ws.Range("D5").Value = "36"
$ws.Range("E5").Value = "36"
$ws.Range("F5").Value = "36"
$ws.Range("G5").Value = "36"

# Update row 6 values: old "36" -> "3" (these are the values that used to live in row 7)
$ws.Range("B6").Value = "3"
$ws.Range("C6").Value = "3"
$ws.Range("D6").Value = "3"
$ws.Range("E6").Value = "3"
$ws.Range("F6").Value = "3"
$ws.Range("G6").Value = "3"

# Restore the normal (unformatted) cell style now that the text is locked in.
$ws.Range("B3").Style = "Normal"
$ws.Range("B5:G5").Style = "Normal"
$ws.Range("B6:G6").Style = "Normal"

# Delete entire row 7 (A7:G7), shrinking the used range down to A1:G6
$ws.Rows.Item(7).Delete()
